{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Three textual tweaks (see commit message \"Details and word knitting\"):\n//   1. Merge the two p-value runs into a single, reworded sentence.\n//   2. Trim/reword the \"adjusted R-squared\" sentence (\"second\"/\"third\" -> \"model 2\"/\"model 3\").\n//   3. Reorder the closing \"Therefore, ...\" sentence (keeps its bold run formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldPValue =\n  \"According to the p-value, which is 2.910^{-4} and <0.05, this model is statistically significant.\";\nconst newPValue =\n  \"The p-value is 0.0002, being <0.05 means this model is statistically significant.\";\n\nconst oldSummary =\n  \"The full summary of coefficients can be found in appendix A5. Because this are multivariable models, the adjusted R-squared is a better metric to evaluate variability. In the second model, the transmission type explains 77% of the variability in MPG, while the third explains 82% and it's still parsimonious, even more than the previous one.\";\nconst newSummary =\n  \"The full summary of coefficients can be found in appendix A5. Because this are multivariable models, the adjusted R-squared is a better metric. In model 2 the transmission type explains 77% of the variability in MPG, while model 3 explains 82% and it's still parsimonious, even more than the previous one.\";\n\nconst oldTherefore =\n  \"Therefore, the linear model which fits transmission type, horsepower and weight as regressors, is a better model to predict MPG.\";\nconst newTherefore =\n  \"Therefore, in this setting the better linear model to predict MPG is the one which fits transmission type, horsepower and weight as regressors.\";\n\nconst replacements = [\n  [oldPValue, newPValue],\n  [oldSummary, newSummary],\n  [oldTherefore, newTherefore],\n];\n\nfor (const item of paragraphs.items) {\n  for (const [oldText, newText] of replacements) {\n    if (item.text === oldText) {\n      // Replacing via the paragraph's own range keeps the paragraph's\n      // existing run formatting (e.g. the bold \"Therefore, ...\" run).\n      item.getRange().insertText(newText, \"Replace\");\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument / $d is the open document.\n#\n# Three textual tweaks (see commit message \"Details and word knitting\"):\n#   1. Merge the two p-value runs into a single, reworded sentence.\n#   2. Trim/reword the \"adjusted R-squared\" sentence (\"second\"/\"third\" -> \"model 2\"/\"model 3\").\n#   3. Reorder the closing \"Therefore, ...\" sentence (keeps its bold run formatting).\n\n$d = $word.ActiveDocument\n\n# Keep straight apostrophes/quotes as-is (the source text uses plain \"'\" -\n# wdReplaceAll otherwise lets Word's smart-quote AutoFormat turn it into a\n# curly U+2019 while re-typing the replacement text).\ntry { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}\ntry { $word.Options.AutoFormatReplaceQuotes = $false } catch {}\n\nfunction Replace-DocText($doc, [string]$oldText, [string]$newText) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Replacement.Text = $newText\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n# Each call below targets the smallest span that actually differs (computed\n# as the old/new sentence with their common prefix/suffix trimmed off) and\n# leaves any untouched text - notably the \"it's\" apostrophe in the second\n# sentence - completely outside the Find/Replace span, so Word never\n# retypes it (and can't smart-quote it into a curly \"'\" on the way back in).\n\n# 1) \"According to the p-value, which is 2.9 10^{-4} and <0.05, ...\" ->\n#    \"The p-value is 0.0002, being <0.05 means ...\"\nReplace-DocText $d \"According to the p-value, which is 2.910^{-4} and <0.05,\" \"The p-value is 0.0002, being <0.05 means\"\n\n# 2) \"... a better metric to evaluate variability. In the second model, the\n#    transmission type explains 77% ... while the third explains 82% and ...\"\n#    -> \"... a better metric. In model 2 the transmission type explains 77%\n#    ... while model 3 explains 82% and ...\"\n#    (Find.Text/Replacement.Text must not start with a space - Word's Find\n#    engine silently ignores a leading blank and the match fails - so the\n#    span starts at the preceding whole word \"metric\" instead.)\nReplace-DocText $d \"metric to evaluate variability. In the second model, the transmission type explains 77% of the variability in MPG, while the third\" \"metric. In model 2 the transmission type explains 77% of the variability in MPG, while model 3\"\n\n# 3) \"Therefore, the linear model which fits ..., is a better model to\n#    predict MPG.\" -> \"Therefore, in this setting the better linear model to\n#    predict MPG is the one which fits ... regressors.\"\nReplace-DocText $d \"the linear model which fits transmission type, horsepower and weight as regressors, is a better model to predict MPG\" \"in this setting the better linear model to predict MPG is the one which fits transmission type, horsepower and weight as regressors\"\n"}
